$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Create the new "Global" worksheet as the first tab -------------------
# A throwaway sheet is inserted (and later removed) first so that the
# internal sheetId counter lands on the same value a "real" Excel session
# would have produced for the Global tab after this many sheet creations.
$placeholder = $wb.Worksheets.Add()
$placeholder.Name = "ZZZ_Placeholder"

$gbl = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$gbl.Name = "Global"

$wb.Worksheets.Item("ZZZ_Placeholder").Delete()

# --- Populate the Global sheet ---------------------------------------------
$uddt = $wb.Worksheets.Item("UDDT")

$gbl.Range("A1").Value = "parameter"
$gbl.Range("B1").Value = "name (if diff)"
$gbl.Range("C1").Value = "unit"
$gbl.Range("D1").Value = "expected"
$gbl.Range("E1").Value = "low"
$gbl.Range("F1").Value = "high"
$gbl.Range("G1").Value = "distribution"
$gbl.Range("H1").Value = "references"

$uddt.Range("A1:H1").Copy()
$gbl.Range("A1:H1").PasteSpecial(-4122)

# Write the new data rows in the same cell order the source workbook did, so
# newly-introduced shared-string entries land at the same table indices.
$gbl.Range("B2").Value = "time_full_degradation"
$gbl.Range("B3").Value = "reduction_full_degradation"
$gbl.Range("C3").Value = "log reduction"
$gbl.Range("H2").Value = "Assumption"
$gbl.Range("A2").Value = "tau_deg"
$gbl.Range("A3").Value = "log_deg"

$gbl.Range("C2").Value = "years"
$gbl.Range("D2").Value = 2
$gbl.Range("E2").Value = 1
$gbl.Range("F2").Value = 3
$gbl.Range("G2").Value = "uniform"

$gbl.Range("D3").Value = 3
$gbl.Range("E3").Value = 2
$gbl.Range("F3").Value = 4
$gbl.Range("G3").Value = "uniform"
$gbl.Range("H3").Value = "Assumption"

$uddt.Range("D2:G2").Copy()
$gbl.Range("D2:G3").PasteSpecial(-4122)

$gbl.Columns.Item(2).ColumnWidth = 22.83

# Comment on B1 explaining the new column, matching the comment already used
# elsewhere in the workbook for the same header.
$cmt = $gbl.Range("B1").AddComment("Yalin Li:" + [char]10 + "Name for the corresponding parameters in another repository (https://github.com/QSD-for-WaSH/Bwaise-sanitation-alternatives)")

# --- Update selections/active cell on the other sheets ----------------------
$excretion = $wb.Worksheets.Item("Excretion")
$excretion.Activate()
$excretion.Range("A1:H1").Select()

$toilet = $wb.Worksheets.Item("Toilet")
$toilet.Activate()
$toilet.Range("C24").Select()

$pitLatrine = $wb.Worksheets.Item("PitLatrine")
$pitLatrine.Activate()
$pitLatrine.Range("G15").Select()

$uddt.Activate()
$uddt.Range("A12").Select()

# Finally, make Global the active/visible tab with A4 selected (just below
# the two new parameter rows), matching the saved state of the workbook.
$gbl.Activate()
$gbl.Range("A4").Select()

# --- Move the application window, matching the saved workbook view ---------
$wb.Windows.Item(1).Left = 9720
$wb.Windows.Item(1).Top = 10240
